# Exemple_EvaluationCouts.xlsx - correction readme clavier et update MAN
#
# - B2: add note "En pack de 5, 150 pour les PCB et 50 pour le transport"
#       in a small (8pt) black font.
# - E2: 5 -> 1 (quantity)
# - F2: 50 -> 10 (shipping)
# - G2 formula (D2*E2)+F2 recalculates automatically.
# - Selection moves to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste")

# Update the PCB row quantity / shipping figures.
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 10

# Add the explanatory note in B2 with a smaller 8pt black font, matching
# the look of the surrounding table (same fill/border, new font).
$ws.Range("B2").Value = "En pack de 5, 150 pour les PCB et 50 pour le transport"
$ws.Range("B2").Font.Size = 8
$ws.Range("B2").Font.Color = 0
$ws.Range("B2").Font.Name = "Calibri"

# Match the selection left behind in the saved workbook.
$ws.Range("B4").Select()

$wb.Save()
